$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Outstanding Balance), shifting it and
# all subsequent liability columns one to the right.
$ws.Columns("C").Insert()

# Populate the new column's header with the new field name.
$ws.Range("C1").Value = "Liability Owner Name"

# Give the new column roughly the same width as columns A and B.
$ws.Columns("C").ColumnWidth = 19.998697916666668

# Move the active selection, matching the saved workbook state.
[void]$ws.Range("D6").Select()
